# Update the dSF column (F) values for specific rows as part of
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 3
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -8
$ws.Range("F13").Value = -5
$ws.Range("F16").Value = -2
